$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 29
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 300
$ws.Range("L2").Value = 83
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 56
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0
$ws.Range("S2").Value = 37
$ws.Range("T2").Value = 63
$ws.Range("V2").Value = 487
$ws.Range("X2").Value = 496
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 7
$ws.Range("AA2").Value = 3
